# Rashid Khan match log - update latest-match activity (row 2) and shift
# the two most recent innings rows (7, 8) down so row 8 now holds the
# newest figures and row 7 the previous (now-cleared) slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 (balls faced in the newest row) flips from 0 to 1
$ws.Range("D2").Value = 1

# Row 7 -> cleared to zeros
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0

# Row 8 -> now carries what row 7 used to have
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 1
